# Generate Report for Handoff
# Updates the generated-file identifiers (old guid -> new guid), the
# handoff xliff file names, and their associated timestamps, across the
# Overview / zh-cn / de-de sheets of the handback status report.

$wb = $excel.ActiveWorkbook

$oldGuid = "6e89b827-02ff-40f7-b03b-1d8db5a11196"
$newGuid = "e71a5e0b-2db9-4931-9633-2e9f660ff036"

$newMdName      = "$newGuid.md"
$newMdPath      = "e2e\$newGuid.md"
$newHoGenDate   = "2016-09-01 11:07:15"
$newHandoffDate = "2016-09-01 11:07:09"
$newZhCnXlf     = "$newGuid.eb7824e3bbd00b5504173680fc256f67544c8060.zh-cn.xlf"
$newDeDeXlf     = "$newGuid.eb7824e3bbd00b5504173680fc256f67544c8060.de-de.xlf"

# ---------------------------------------------------------------------
# Sheet "Overview"
#   A2 File Name, B2 Path And Name, G2 Latest HO Xliff Generate Date
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newMdName
$wsOverview.Range("B2").Value = $newMdPath
$wsOverview.Range("G2").Value = $newHoGenDate

foreach ($link in $wsOverview.Hyperlinks) {
    $link.TextToDisplay = $newMdPath
}

# ---------------------------------------------------------------------
# Sheet "zh-cn"
#   A2 Source File Name, G2 Latest Handoff File, H2 Latest Handoff Datetime
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = $newMdName
$wsZhCn.Range("G2").Value = $newZhCnXlf
$wsZhCn.Range("H2").Value = $newHandoffDate

foreach ($link in $wsZhCn.Hyperlinks) {
    $link.TextToDisplay = $newMdName
}

# ---------------------------------------------------------------------
# Sheet "de-de"
#   A2 Source File Name, G2 Latest Handoff File, H2 Latest Handoff Datetime
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = $newMdName
$wsDeDe.Range("G2").Value = $newDeDeXlf
$wsDeDe.Range("H2").Value = $newHoGenDate

foreach ($link in $wsDeDe.Hyperlinks) {
    $link.TextToDisplay = $newMdName
}
